# Share of Cargo Dist Transported that is New This Year
# Split the single "SoCDTtiNTY" sheet (passenger vs. freight columns) into
# two sheets - one per cargo type - each broken out by vehicle
# (powertrain) type across columns B:H.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoCDTtiNTY")

# --- capture the original passenger (col B) / freight (col C) values ----
# (rows 2-7) before we start overwriting anything.
$psgrVals = @()
$frgtVals = @()
for ($r = 2; $r -le 7; $r++) {
    $psgrVals += $ws.Cells.Item($r, 2).Value2
    $frgtVals += $ws.Cells.Item($r, 3).Value2
}

$vehicleTypes = @(
    "battery electric vehicle",
    "natural gas vehicle",
    "gasoline vehicle",
    "diesel vehicle",
    "plugin hybrid vehicle",
    "LPG vehicle",
    "hydrogen vehicle"
)

# --- rename the existing sheet to the "passenger" variant ---------------
$ws.Name = "SoCDTtiNTY-psgr"

# --- duplicate it to create the "freight" variant ------------------------
$ws.Copy($null, $ws)
$wsFrgt = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsFrgt.Name = "SoCDTtiNTY-frgt"

function Fill-Sheet($sheet, $values) {
    # column widths
    $sheet.Columns.Item(1).ColumnWidth = 19.140625
    $sheet.Range($sheet.Columns.Item(2), $sheet.Columns.Item(8)).ColumnWidth = 14.42578125

    # header row
    $sheet.Range("A1").Value = "Share that is New (dimensionless)"
    $sheet.Range("A1").Font.Bold = $true
    $sheet.Range("A1").WrapText = $true

    for ($c = 2; $c -le 8; $c++) {
        $cell = $sheet.Cells.Item(1, $c)
        $cell.Value = $vehicleTypes[$c - 2]
        $cell.Font.Bold = $false
        $cell.WrapText = $true
        $cell.HorizontalAlignment = -4152   # xlRight
    }
    $sheet.Rows.Item(1).RowHeight = 30

    # data rows: broadcast the captured value across columns B:H
    for ($i = 0; $i -lt 6; $i++) {
        $r = $i + 2
        for ($c = 2; $c -le 8; $c++) {
            $sheet.Cells.Item($r, $c).Value = $values[$i]
        }
    }
}

Fill-Sheet $ws $psgrVals
Fill-Sheet $wsFrgt $frgtVals
